$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.504.70"
$ws.Range("E2").Value = "  -5.71%  "
$ws.Range("D3").Value = "3.283.49"
$ws.Range("E3").Value = "  -6.04%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "550.95"
$ws.Range("E5").Value = "  -4.50%  "
$ws.Range("D6").Value = "182.45"
$ws.Range("E6").Value = "  -3.95%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.580"
$ws.Range("E8").Value = "  -5.00%  "
$ws.Range("D9").Value = "3.276.02"
$ws.Range("E9").Value = "  -5.98%  "
$ws.Range("D10").Value = "0.179"
$ws.Range("E10").Value = "  -12.14%  "
$ws.Range("D11").Value = "0.578"
$ws.Range("E11").Value = "  -6.43%  "
$ws.Range("D12").Value = "46.98"
$ws.Range("E12").Value = "  -7.17%  "
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -7.30%  "
$ws.Range("D14").Value = "8.61"
$ws.Range("E14").Value = "  -5.36%  "
$ws.Range("D15").Value = "629.98"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "3.796.41"
$ws.Range("E16").Value = "  -6.49%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "18.02"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "65.475.77"
$ws.Range("E18").Value = "  -5.58%  "
$ws.Range("E19").Value = "  -3.30%  "
$ws.Range("D20").Value = "3.266.90"
$ws.Range("E20").Value = "  -6.45%  "
$ws.Range("D21").Value = "11.28"
$ws.Range("E21").Value = "  -8.24%  "
$ws.Range("D22").Value = "0.898"
$ws.Range("E22").Value = "  -5.19%  "
$ws.Range("D23").Value = "17.67"
$ws.Range("E23").Value = "  -3.38%  "
$ws.Range("D24").Value = "104.75"
$ws.Range("E24").Value = "  +6.20%  "
$ws.Range("D25").Value = "4.93"
$ws.Range("E25").Value = "  -6.62%  "
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -9.33%  "
$ws.Range("D27").Value = "2.66"
$ws.Range("E27").Value = "  -7.24%  "
$ws.Range("D28").Value = "9.44"
$ws.Range("E28").Value = "  -4.84%  "
$ws.Range("D29").Value = "8.57"
$ws.Range("E29").Value = "  -7.87%  "
$ws.Range("D30").Value = "29.87"
$ws.Range("E30").Value = "  -7.65%  "
$ws.Range("D31").Value = "6.28"
$ws.Range("E31").Value = "  -5.57%  "
$ws.Range("D32").Value = "3.75"
$ws.Range("E32").Value = "  -7.60%  "
$ws.Range("D33").Value = "11.01"
$ws.Range("E33").Value = "  -4.41%  "
$ws.Range("D34").Value = "0.104"
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.765.59"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "530.42"
$ws.Range("E36").Value = "  -10.66%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "56.90"
$ws.Range("E38").Value = "  -6.45%  "
$ws.Range("D39").Value = "0.0₃0727"
$ws.Range("E39").Value = "  -7.69%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "34.01"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "3.24"
$ws.Range("E41").Value = "  -8.40%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.128"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  -6.24%  "
$ws.Range("D44").Value = "3.23"
$ws.Range("E44").Value = "  -14.26%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.330"
$ws.Range("E45").Value = "  -10.85%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0411"
$ws.Range("E46").Value = "  -6.33%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.20"
$ws.Range("E47").Value = "  -3.33%  "
$ws.Range("D48").Value = "0.128"
$ws.Range("E48").Value = "  -5.17%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").Value = "2.55"
$ws.Range("E49").Value = "  -9.63%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "1.24"
$ws.Range("E51").Value = "  +1.57%  "
